$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 is alinea "1f" ("Desenvolvimento de salvar imagem cropada"):
# mark it as complete -> progress 100%, owner Bernardo, status "Done!".
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = "Bernardo"

# E13 used to be the master cell of the shared "-" formula covering
# E12:E27 (column "Encarregue" shows "-" while there's no progress).
# Re-apply that formula contiguously over E14:E27 so the group's new
# master (E14) keeps carrying it, then restore the handful of rows in
# that span that actually hold a person's name rather than "-".
$ws.Range("E14:E27").Formula = '="-"'
$ws.Range("E15").Value = "João"
$ws.Range("E22").Value = "Eduardo"
$ws.Range("E23").Value = "Eduardo"
$ws.Range("E25").Value = "João"

# Update the active selection to H14, as in the saved workbook.
$ws.Range("H14").Select()
